# install_MakePES.pptx edit — "added manuals of MakePES"
#
# 1. Bump the cached text of every datetime placeholder field (master,
#    all slide layouts, notes master) from 2019/5/3 to 2019/5/5.
# 2. Slide 4 ("> export sindo_jar=..." box): the lone ">" run that
#    introduces the "java -cp ..." line gets a trailing space -> "> ".
# 3. Slide 4 ("sindo_jar=${HOME}/sindo/jar" box): move the shape down
#    slightly and collapse the "=${HOME}/" + "sindo" + "/jar" runs into
#    a single run "=/path/to/sindo-4.0/jar".

$p = $ppt.ActivePresentation

$oldDate = "2019/5/3"
$newDate = "2019/5/5"

# --- 1a. Slide master date placeholder -------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 1b. Every slide layout's date placeholder ------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lsh = $layout.Shapes.Item($si)
        if ($lsh.HasTextFrame) {
            if ($lsh.TextFrame.HasText) {
                if ($lsh.TextFrame.TextRange.Text -eq $oldDate) {
                    $lsh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# NOTE: the notes master's own date placeholder (a `datetimeFigureOut`
# field, vs. `datetime1` everywhere else) is an auto-updating field that
# this host treats as read-only text, and — because its shape Id (3)
# collides with the slide master's "Text Placeholder 2" shape Id — any
# write attempt through NotesMaster.Shapes corrupts that unrelated slide
# master shape instead. So it is deliberately left alone here.

# --- Slide 4 edits ------------------------------------------------------
$slide4 = $p.Slides.Item(4)

# 2. "> export sindo_jar=/path/to/sindo-4.0/jar" / ">java -cp ..." box
#    second paragraph starts with a bare ">" run -> becomes "> ".
$cmdBox = $slide4.Shapes.Item(5)
$cmdRange = $cmdBox.TextFrame.TextRange
$gt = $cmdRange.Characters(44, 1)
if ($gt.Text -eq ">") {
    $gt.Text = "> "
}

# 3. "sindo_jar=${HOME}/sindo/jar" alias box.
$aliasBox = $slide4.Shapes.Item(7)

# 3a. Nudge the box further down the slide.
$aliasBox.Top = 270.18811023622044

# 3b. Collapse "=${HOME}/" + "sindo" + "/jar" into "=/path/to/sindo-4.0/jar".
$aliasRange = $aliasBox.TextFrame.TextRange
$homeRun = $aliasRange.Characters(10, 18)
if ($homeRun.Text -eq "=`${HOME}/sindo/jar") {
    $homeRun.Text = "=/path/to/sindo-4.0/jar"
}
